# Update countries & provincias Spain
# - Refresh "last updated" timestamp
# - Refresh COVID case numbers for a handful of countries whose case counts
#   caused them to change rank (Bolivia overtakes Barein, Honduras overtakes
#   Chequia/Camerun, Islas Malvinas/Santa Sede move ahead of
#   Groenlandia/Islas Turcas y Caicos), plus a data refresh for Nueva Zelanda.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Title / timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 18 de Junio de 2020 a las 05:05"

# Helper to write one full data row:
#   País, Casos totales, Nuevos casos, Casos activos, Recuperados,
#   Casos criticos, Muertes hoy, Muertes
function Set-CountryRow($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Bolivia now ranks above Barein
Set-CountryRow 49 "Bolivia" 20685 802 4002 16004 0 20 679
Set-CountryRow 50 "Barein"  19961 0   14185 5727 0 0  49

# Honduras now ranks above Chequia and Camerun
Set-CountryRow 65 "Honduras" 10299 643 1116 8847 0 6 336
Set-CountryRow 66 "Chequia"  10162 0   7399 2430 0 0 333
Set-CountryRow 67 "Camerun"  9864  0   5570 4018 0 0 276

# Data refresh for Nueva Zelanda (no rank change)
Set-CountryRow 115 "Nueva Zelanda" 1507 1 1482 3 0 0 22

# Islas Malvinas now ranks above Groenlandia, and Santa Sede above Islas
# Turcas y Caicos
Set-CountryRow 206 "Islas Malvinas" 13 0 13 0 0 0 0
Set-CountryRow 207 "Groenlandia"    13 0 13 0 0 0 0
Set-CountryRow 208 "Santa Sede"               12 0 12 0 0 0 0
Set-CountryRow 209 "Islas Turcas y Caicos"    12 0 11 0 0 0 1
